$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.925.11'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +2.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.342.30'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.996'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.34'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +5.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.97'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.359.78'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('E10').Value = '  +1.87%  '
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('E12').Value = '  +1.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.356'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +5.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.780.14'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.54'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.968.39'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.348.64'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.58'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.25'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +2.81%  '
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.74'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.87'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('E25').Value = '  +4.37%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.52'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.44'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +11.14%  '
$ws.Range('E29').Value = '  +4.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.34'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('E31').Value = '  +2.66%  '
$ws.Range('E32').Value = '  +0.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.57'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('E34').Value = '  +15.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.995'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.27'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.16'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +5.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.64'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +5.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '39.38'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +2.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '150.47'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('E42').Value = '  +1.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.64'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '284.62'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.29'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +6.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0931'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0506'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +2.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.563'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('E49').Value = '  +2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.64'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +2.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.381'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +0.36%  '
